# Regenerate column G ("K") values for gsellman_robert.xlsx save_data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, keyed by row number (row 23 is unchanged and intentionally omitted).
$newK = [ordered]@{
    2  = 1
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 2
    9  = 3
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 2
    21 = 1
    22 = 0
    24 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
